# Weekly data refresh: add the 2021-01-08 (serial 44204) week to both the
# national Fallecido_Recuperado summary sheet and the per-province
# Provincias_Semanal sheet, following the existing pattern of "header date
# row, then filled-in-next-week" rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Provincias_Semanal: complete the pending 2021-01-01 (44197) week block
# (row 1282 onward) with all 32 provinces, then append the new week's
# leading date row (1314) for 2021-01-08 (44204).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Provincias_Semanal")

$ws2.Range("B1282").Value = "Distrito Nacional"
$ws2.Range("C1282").Value = 4592.57
$ws2.Range("D1282").Value = 412

$ws2.Range("A1283").Value = 44197
$ws2.Range("B1283").Value = "Azua"
$ws2.Range("C1283").Value = 996.6
$ws2.Range("D1283").Value = 26

$ws2.Range("A1284").Value = 44197
$ws2.Range("B1284").Value = "Baoruco"
$ws2.Range("C1284").Value = 1053.4100000000001
$ws2.Range("D1284").Value = 7

$ws2.Range("A1285").Value = 44197
$ws2.Range("B1285").Value = "Barahona"
$ws2.Range("C1285").Value = 1188.31
$ws2.Range("D1285").Value = 23

$ws2.Range("A1286").Value = 44197
$ws2.Range("B1286").Value = "Dajabon"
$ws2.Range("C1286").Value = 842.52
$ws2.Range("D1286").Value = 7

$ws2.Range("A1287").Value = 44197
$ws2.Range("B1287").Value = "Duarte"
$ws2.Range("C1287").Value = 1720.9
$ws2.Range("D1287").Value = 125

$ws2.Range("A1288").Value = 44197
$ws2.Range("B1288").Value = "Elias Pina"
$ws2.Range("C1288").Value = 425.85
$ws2.Range("D1288").Value = 5

$ws2.Range("A1289").Value = 44197
$ws2.Range("B1289").Value = "El Seibo"
$ws2.Range("C1289").Value = 638.45000000000005
$ws2.Range("D1289").Value = 7

$ws2.Range("A1290").Value = 44197
$ws2.Range("B1290").Value = "Espaillat"
$ws2.Range("C1290").Value = 1362.55
$ws2.Range("D1290").Value = 89

$ws2.Range("A1291").Value = 44197
$ws2.Range("B1291").Value = "Independencia"
$ws2.Range("C1291").Value = 1371.01
$ws2.Range("D1291").Value = 9

$ws2.Range("A1292").Value = 44197
$ws2.Range("B1292").Value = "La Altagracia"
$ws2.Range("C1292").Value = 1778.12
$ws2.Range("D1292").Value = 44

$ws2.Range("A1293").Value = 44197
$ws2.Range("B1293").Value = "La Romana"
$ws2.Range("C1293").Value = 1950.5
$ws2.Range("D1293").Value = 112

$ws2.Range("A1294").Value = 44197
$ws2.Range("B1294").Value = "La Vega"
$ws2.Range("C1294").Value = 1699.29
$ws2.Range("D1294").Value = 126

$ws2.Range("A1295").Value = 44197
$ws2.Range("B1295").Value = "Maria Trinidad Sanchez"
$ws2.Range("C1295").Value = 1265.8
$ws2.Range("D1295").Value = 13

$ws2.Range("A1296").Value = 44197
$ws2.Range("B1296").Value = "Monte Cristi"
$ws2.Range("C1296").Value = 607.4
$ws2.Range("D1296").Value = 14

$ws2.Range("A1297").Value = 44197
$ws2.Range("B1297").Value = "Pedernales"
$ws2.Range("C1297").Value = 1637.28
$ws2.Range("D1297").Value = 3

$ws2.Range("A1298").Value = 44197
$ws2.Range("B1298").Value = "Peravia"
$ws2.Range("C1298").Value = 793.18
$ws2.Range("D1298").Value = 44

$ws2.Range("A1299").Value = 44197
$ws2.Range("B1299").Value = "Puerto Plata"
$ws2.Range("C1299").Value = 1532.02
$ws2.Range("D1299").Value = 131

$ws2.Range("A1300").Value = 44197
$ws2.Range("B1300").Value = "Hermanas Mirabal"
$ws2.Range("C1300").Value = 1379.3
$ws2.Range("D1300").Value = 25

$ws2.Range("A1301").Value = 44197
$ws2.Range("B1301").Value = "Samana"
$ws2.Range("C1301").Value = 556.35
$ws2.Range("D1301").Value = 3

$ws2.Range("A1302").Value = 44197
$ws2.Range("B1302").Value = "San Cristobal"
$ws2.Range("C1302").Value = 878.84
$ws2.Range("D1302").Value = 119

$ws2.Range("A1303").Value = 44197
$ws2.Range("B1303").Value = "San Juan"
$ws2.Range("C1303").Value = 1177.53
$ws2.Range("D1303").Value = 45

$ws2.Range("A1304").Value = 44197
$ws2.Range("B1304").Value = "San Pedro de Macoris"
$ws2.Range("C1304").Value = 702.05
$ws2.Range("D1304").Value = 37

$ws2.Range("A1305").Value = 44197
$ws2.Range("B1305").Value = "Sanchez Ramirez"
$ws2.Range("C1305").Value = 1606.81
$ws2.Range("D1305").Value = 18

$ws2.Range("A1306").Value = 44197
$ws2.Range("B1306").Value = "Santiago"
$ws2.Range("C1306").Value = 1723.74
$ws2.Range("D1306").Value = 359

$ws2.Range("A1307").Value = 44197
$ws2.Range("B1307").Value = "Santiago Rodriguez"
$ws2.Range("C1307").Value = 1381.27
$ws2.Range("D1307").Value = 10

$ws2.Range("A1308").Value = 44197
$ws2.Range("B1308").Value = "Valverde"
$ws2.Range("C1308").Value = 718.09
$ws2.Range("D1308").Value = 28

$ws2.Range("A1309").Value = 44197
$ws2.Range("B1309").Value = "Monsenor Nouel"
$ws2.Range("C1309").Value = 1583.68
$ws2.Range("D1309").Value = 34

$ws2.Range("A1310").Value = 44197
$ws2.Range("B1310").Value = "Monte Plata"
$ws2.Range("C1310").Value = 386.32
$ws2.Range("D1310").Value = 26

$ws2.Range("A1311").Value = 44197
$ws2.Range("B1311").Value = "Hato Mayor"
$ws2.Range("C1311").Value = 642.48
$ws2.Range("D1311").Value = 12

$ws2.Range("A1312").Value = 44197
$ws2.Range("B1312").Value = "San Jose de Ocoa"
$ws2.Range("C1312").Value = 1231.5
$ws2.Range("D1312").Value = 12

$ws2.Range("A1313").Value = 44197
$ws2.Range("B1313").Value = "Santo Domingo"
$ws2.Range("C1313").Value = 1199.3399999999999
$ws2.Range("D1313").Value = 491

# New week's leading date row -- copy A1282's date formatting (numFmtId 14)
# onto the new cell before writing the value, same as the source sheet.
$ws2.Range("A1282").Copy()
$ws2.Range("A1314").PasteSpecial(-4122)
$ws2.Range("A1314").Value = 44204

$ws2.Range("B1314").Select()

# ---------------------------------------------------------------------
# Fallecido_Recuperado: fill in the pending 44197 week row (42), then
# append the new week's leading date row (43) for 44204. This sheet is
# edited/selected last so it ends up the active tab, matching the source.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Fallecido_Recuperado")

$ws1.Range("B42").Value = 172965
$ws1.Range("C42").Value = 2416
$ws1.Range("D42").Value = 132282

$ws1.Range("A42").Copy()
$ws1.Range("A43").PasteSpecial(-4122)
$ws1.Range("A43").Value = 44204

$ws1.Range("B43").Select()
